$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# H column values (col 8) and I column formulas (col 9)
$hValues = @(2,2,3,2,2,1,1,1,2,3,2,2,2,1,3,1,3,2,2,1,1,2,2,2,3,2,1,1,2,1,2,2,3,3,2,2,2,2,1,3,2,3,2,1,1,2,3,2,2,2,2,2,1,2,2,2,2,3,2,2,2,1,2,2,1,2,2,2,1,1,2,2,2,2,2,1,2,2,2,2,1,3,2,2,2,2,1,1,1,2,2,2,1,2,1,1,2,2,2,2,1,2,1,2,2,1,1,2,1,2,2,2,2,2,2,1,1,2,1,2,2,2)
$bHighlight = @(1,1,1,1,1,0,1,0,0,1,0,1,1,1,0,1,0,1,1,0,1,1,1,1,0,0,1,1,0,1,0,1,0,0,1,1,1,1,1,1,1,0,0,0,0,0,1,0,1,1,1,0,0,1,1,1,1,0,1,1,0,0,0,1,1,1,0,1,0,0,1,1,1,1,1,0,0,1,1,1,0,0,0,0,1,0,1,1,0,0,1,0,0,1,0,1,1,0,1,1,1,1,0,0,1,1,1,1,0,0,1,1,1,1,0,1,1,1,0,1,1,1)
$rowNums = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58,59,60,61,62,63,64,65,66,67,68,69,70,71,72,73,74,75,76,77,78,79,80,81,82,83,84,85,86,87,88,89,90,91,92,93,94,95,96,97,98,99,100,101,102,103,104,105,106,107,108,109,110,111,112,113,114,115,116,117,118,119,120,121,122,123)

for ($idx = 0; $idx -lt $rowNums.Length; $idx++) {
    $r = $rowNums[$idx]
    $ws.Cells.Item($r, 8).Value = $hValues[$idx]
    $ws.Cells.Item($r, 8).NumberFormat = $ws.Cells.Item($r, 4).NumberFormat
    $ws.Cells.Item($r, 9).FormulaR1C1 = "=GEOMEAN(RC[-5]:RC[-1])"
    if ($bHighlight[$idx] -eq 1) {
        $ws.Cells.Item($r, 2).Interior.Color = 65535
    }
}

# Sheet view changes
$ws.Range("C6").Select()
